$d = $word.ActiveDocument

# --- 1. "Sprite :" paragraph: add a second tab before "spr_<Name>" ---
$d.Content.Find.Execute("Sprite :`tspr_", $false, $false, $false, $false, $false, `
                        $true, 1, $false, "Sprite :`t`tspr_", 2)

# --- 2. Insert a brand new "Background:`tbg_<Name>" paragraph right after "Sprite :" ---
$spritePara = $d.Paragraphs.Item(10)
$spritePara.Range.InsertParagraphAfter()
$bgPara = $d.Paragraphs.Item(11)
$bgPara.Range.Text = "Background:`tbg_<Name>"

# --- 3. "Son :" paragraph: drop the en-US language override and add a second tab ---
$sonPara = $d.Paragraphs.Item(12)
$sonPara.Range.Delete()
$bgParaAfter = $d.Paragraphs.Item(11)
$bgParaAfter.Range.InsertParagraphAfter()
$newSonPara = $d.Paragraphs.Item(12)
$newSonPara.Range.Text = "Son :`t`tsnd_<Name>"

Write-Host "Done"
